# set.ambiguous.to is NA by default
#
# Rename the raw presence-records sheet from "M. natalensis" to
# "Raw_M_natalensis_presences" and make it the active tab (it was
# previously the Cleaned_M_natalensis_presences sheet that was active).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("M. natalensis")

# Renaming cascades into every _xlnm.Print_Titles_* duplicate + the wvu
# custom-view names that reference the sheet, but the sheet's own
# PrintArea / PrintTitleRows (backing _xlnm.Print_Area / _xlnm.Print_Titles)
# need to be re-applied explicitly so they pick up the new sheet name too.
$ws.Name = "Raw_M_natalensis_presences"
$ws.PageSetup.PrintArea = "A1:J695"
$ws.PageSetup.PrintTitleRows = "$1:$1"

# Make this newly-renamed sheet the active/selected tab instead of
# Cleaned_M_natalensis_presences.
$ws.Activate()
